# lesson-04.pptx update:
#  1. Insert a new slide ("Правила хорошего тона") as slide #10, built by
#     duplicating slide #9 (same title-rectangle template) and replacing its
#     body content with a fresh text box.
#  2. Bump the handout/footer date field on every slide/layout/master from
#     8/17/2012 to 8/21/2012.

$p = $ppt.ActivePresentation

# --- 1. New slide -----------------------------------------------------
$src = $p.Slides.Item(9)
$new = $src.Duplicate()
$ns = $p.Slides.Item(10)

# Drop everything but the title rectangle (shape 1); the source slide
# carries 7 shapes, the new slide only needs the title + one text box.
for ($i = $ns.Shapes.Count; $i -ge 2; $i--) {
    $ns.Shapes.Item($i).Delete()
}

$title = $ns.Shapes.Item(1)
$title.TextFrame.TextRange.Text = "Правила хорошего тона"
$title.Left = 204.37503937007875
$title.Top = 0
$title.Width = 265.27055118110235
$title.Height = 36.3515748031496

$body = $ns.Shapes.AddTextbox(1, 25.474645669291338, 88.56251968503938, 669.0507086614173, 181.75779527559055)
$body.TextFrame.WordWrap = -1
$body.TextFrame.AutoSize = 1
$body.TextFrame.TextRange.Text = "Windows – многозадачная, многопользовательская, сетевая ОС и программы должны это учитывать.`r`rРеестр (registry) – HKLM, HKCU`r`rПрофили пользователя. Роуминг (roaming) профиля`r`r%Program Files% - только для чтения!"
$body.Left = 25.474645669291338
$body.Top = 88.56251968503938
$body.Width = 669.0507086614173
$body.Height = 181.75779527559055

# --- 2. Footer date field, every slide ---------------------------------
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $sl = $p.Slides.Item($i)
    for ($j = 1; $j -le $sl.Shapes.Count; $j++) {
        $shp = $sl.Shapes.Item($j)
        if ($shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "8/17/2012") {
                $tr.Text = "8/21/2012"
            }
        }
    }
}
